$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.597803197262351
$ws.Range("C2").Value = 0.7119485641435404
$ws.Range("D2").Value = 0.03291285572721492
$ws.Range("E2").Value = 1.361407321067915
$ws.Range("F2").Value = 0.6798212621317248
$ws.Range("G2").Value = 0.0007871586232472039
$ws.Range("H2").Value = 0.01199277220940109
$ws.Range("I2").Value = 0.005767783202876409
$ws.Range("P2").Value = 0.7161795053628026
$ws.Range("Q2").Value = 2.014744168414211

$ws.Range("B3").Value = 2.26066291172782
$ws.Range("C3").Value = 0.6281112901472454
$ws.Range("D3").Value = 0.02978638338481687
$ws.Range("E3").Value = 1.18640414660436
$ws.Range("F3").Value = 0.6162561042387509
$ws.Range("G3").Value = 0.0007907025278933496
$ws.Range("H3").Value = 0.008635475661823488
$ws.Range("I3").Value = 0.003802278604460696
$ws.Range("P3").Value = 0.7252262043596858
$ws.Range("Q3").Value = 1.84802363948063

$ws.Range("B4").Value = 2.053033621404154
$ws.Range("C4").Value = 0.5769465878343851
$ws.Range("D4").Value = 0.02786294938359291
$ws.Range("E4").Value = 1.079243740448391
$ws.Range("F4").Value = 0.5777840511071304
$ws.Range("G4").Value = 0.0007929509702298748
$ws.Range("H4").Value = 0.006788786149038328
$ws.Range("I4").Value = 0.002808140480463628
$ws.Range("P4").Value = 0.7314904468061201
$ws.Range("Q4").Value = 1.747272484144673

$ws.Range("B5").Value = 1.966867226263588
$ws.Range("C5").Value = 0.5570060540242423
$ws.Range("D5").Value = 0.02711732349066409
$ws.Range("E5").Value = 1.03560817179681
$ws.Range("F5").Value = 0.561636489855097
$ws.Range("G5").Value = 0.0007938913891312019
$ws.Range("H5").Value = 0.006084113324819063
$ws.Range("I5").Value = 0.002507873075721534
$ws.Range("P5").Value = 0.734578658242711
$ws.Range("Q5").Value = 1.704511633347835

$ws.Range("B6").Value = 1.950852273096388
$ws.Range("C6").Value = 0.5547293813975784
$ws.Range("D6").Value = 0.02704163866897247
$ws.Range("E6").Value = 1.028341575679136
$ws.Range("F6").Value = 0.5582333237960952
$ws.Range("G6").Value = 0.0007940555458227334
$ws.Range("H6").Value = 0.005967871774951172
$ws.Range("I6").Value = 0.00253600607557658
$ws.Range("P6").Value = 0.7355487681603776
$ws.Range("Q6").Value = 1.694901505474832

$ws.Range("B7").Value = 2.047220876175572
$ws.Range("C7").Value = 0.579499765070608
$ws.Range("D7").Value = 0.02798496047145704
$ws.Range("E7").Value = 1.078590135761104
$ws.Range("F7").Value = 0.5755708414071776
$ws.Range("G7").Value = 0.0007929820356354663
$ws.Range("H7").Value = 0.006772998661915941
$ws.Range("I7").Value = 0.002993613051261335
$ws.Range("P7").Value = 0.7327682632237398
$ws.Range("Q7").Value = 1.739771223154349

$ws.Range("B8").Value = 2.475449915545028
$ws.Range("C8").Value = 0.6867448461070467
$ws.Range("D8").Value = 0.03201250899027031
$ws.Range("E8").Value = 1.300903148091805
$ws.Range("F8").Value = 0.655126001160383
$ws.Range("G8").Value = 0.0007883887738010078
$ws.Range("H8").Value = 0.01077826887721478
$ws.Range("I8").Value = 0.00526474178485703
$ws.Range("P8").Value = 0.7208571998837741
$ws.Range("Q8").Value = 1.947695453102057

$ws.Range("B9").Value = 3.31937571960276
$ws.Range("C9").Value = 0.89486404983694
$ws.Range("D9").Value = 0.03964780683899249
$ws.Range("E9").Value = 1.74031629218878
$ws.Range("F9").Value = 0.819821406831295
$ws.Range("G9").Value = 0.0007799115522071909
$ws.Range("H9").Value = 0.02055750911132859
$ws.Range("I9").Value = 0.01132324856504763
$ws.Range("P9").Value = 0.7010876780096567
$ws.Range("Q9").Value = 2.382696616817753

$ws.Range("B10").Value = 3.917712953324951
$ws.Range("C10").Value = 1.045199349042434
$ws.Range("D10").Value = 0.04618919566138402
$ws.Range("E10").Value = 1.963957886527425
$ws.Range("F10").Value = 0.93232782730243
$ws.Range("G10").Value = 0.0007741951350906243
$ws.Range("H10").Value = 0.02843873462482227
$ws.Range("I10").Value = 0.01710822360147279
$ws.Range("P10").Value = 0.696941688860889
$ws.Range("Q10").Value = 2.671041276445521

$ws.Range("B11").Value = 4.026737869601959
$ws.Range("C11").Value = 1.062850731216088
$ws.Range("D11").Value = 0.05755877076850879
$ws.Range("E11").Value = 1.273281716279527
$ws.Range("F11").Value = 0.8817242950077855
$ws.Range("G11").Value = 0.0007731507858902821
$ws.Range("H11").Value = 0.04364912973931112
$ws.Range("I11").Value = 0.01881957031504289
$ws.Range("P11").Value = 0.7518756670614692
$ws.Range("Q11").Value = 2.455349917182559

$ws.Range("B12").Value = 3.998808706553689
$ws.Range("C12").Value = 1.04164719867569
$ws.Range("D12").Value = 0.06634448787927028
$ws.Range("E12").Value = 0.7757246460136571
$ws.Range("F12").Value = 0.8193973257164942
$ws.Range("G12").Value = 0.0007732720225684575
$ws.Range("H12").Value = 0.07933126607416341
$ws.Range("I12").Value = 0.01876606120865443
$ws.Range("P12").Value = 0.8018015336147926
$ws.Range("Q12").Value = 2.230018738977265

$ws.Range("B13").Value = 3.856906087453297
$ws.Range("C13").Value = 0.9941221507705222
$ws.Range("D13").Value = 0.07374033173330474
$ws.Range("E13").Value = 0.3981850151505242
$ws.Range("F13").Value = 0.7415224755416716
$ws.Range("G13").Value = 0.0007743164300377484
$ws.Range("H13").Value = 0.1320877733521399
$ws.Range("I13").Value = 0.01762405812788259
$ws.Range("P13").Value = 0.8516298066913066
$ws.Range("Q13").Value = 1.974162583525896

$ws.Range("B14").Value = 3.705445098923747
$ws.Range("C14").Value = 0.9502030679651909
$ws.Range("D14").Value = 0.07837039433519521
$ws.Range("E14").Value = 0.2067866589870491
$ws.Range("F14").Value = 0.6801878403898911
$ws.Range("G14").Value = 0.0007754502749329242
$ws.Range("H14").Value = 0.1794594729324785
$ws.Range("I14").Value = 0.01645342251417592
$ws.Range("P14").Value = 0.8871022934417567
$ws.Range("Q14").Value = 1.782298678331358

$ws.Range("B15").Value = 3.642840219142954
$ws.Range("C15").Value = 0.9350401674155364
$ws.Range("D15").Value = 0.07915881600541752
$ws.Range("E15").Value = 0.1691515420980636
$ws.Range("F15").Value = 0.6612674466339854
$ws.Range("G15").Value = 0.0007759511119717162
$ws.Range("H15").Value = 0.1913068405711869
$ws.Range("I15").Value = 0.01602642478418481
$ws.Range("P15").Value = 0.8957326648395139
$ws.Range("Q15").Value = 1.726423683432074

$ws.Range("B16").Value = 3.415509897697746
$ws.Range("C16").Value = 0.8831420113135096
$ws.Range("D16").Value = 0.0746227652443352
$ws.Range("E16").Value = 0.1649757372158618
$ws.Range("F16").Value = 0.6290727610731253
$ws.Range("G16").Value = 0.000778154564283884
$ws.Range("H16").Value = 0.1761469891596477
$ws.Range("I16").Value = 0.01389576357840028
$ws.Range("P16").Value = 0.8846610237793584
$ws.Range("Q16").Value = 1.655776269964292

$ws.Range("B17").Value = 3.321024665012828
$ws.Range("C17").Value = 0.8665758524583111
$ws.Range("D17").Value = 0.068686150943293
$ws.Range("E17").Value = 0.2495924249046269
$ws.Range("F17").Value = 0.6365869619460085
$ws.Range("G17").Value = 0.0007792633737742622
$ws.Range("H17").Value = 0.1374730480102357
$ws.Range("I17").Value = 0.01290158642736738
$ws.Range("P17").Value = 0.8577170160333907
$ws.Range("Q17").Value = 1.70188707370113

$ws.Range("B18").Value = 3.333655881952268
$ws.Range("C18").Value = 0.8764144046545539
$ws.Range("D18").Value = 0.06103533650605186
$ws.Range("E18").Value = 0.4836990774361354
$ws.Range("F18").Value = 0.6812400819957105
$ws.Range("G18").Value = 0.000779451795486122
$ws.Range("H18").Value = 0.08540656275748404
$ws.Range("I18").Value = 0.0125143417278295
$ws.Range("P18").Value = 0.8143454827812491
$ws.Range("Q18").Value = 1.86245199742703

$ws.Range("B19").Value = 3.423425903977261
$ws.Range("C19").Value = 0.9129069093284556
$ws.Range("D19").Value = 0.05334773736532838
$ws.Range("E19").Value = 0.9192880719533036
$ws.Range("F19").Value = 0.7509527051255276
$ws.Range("G19").Value = 0.0007787894674839167
$ws.Range("H19").Value = 0.04300193645821793
$ws.Range("I19").Value = 0.01316157143171637
$ws.Range("P19").Value = 0.7672507512976097
$ws.Range("Q19").Value = 2.100504158686675

$ws.Range("B20").Value = 3.744961384763428
$ws.Range("C20").Value = 1.014395815254744
$ws.Range("D20").Value = 0.04493882850729136
$ws.Range("E20").Value = 1.900556197826361
$ws.Range("F20").Value = 0.8957133485068596
$ws.Range("G20").Value = 0.0007757310825716802
$ws.Range("H20").Value = 0.02618642432198026
$ws.Range("I20").Value = 0.01606431092290084
$ws.Range("P20").Value = 0.7024029010301973
$ws.Range("Q20").Value = 2.571054206805428

$ws.Range("B21").Value = 4.234385973037206
$ws.Range("C21").Value = 1.13884384872128
$ws.Range("D21").Value = 0.04853039978854667
$ws.Range("E21").Value = 2.232452845818756
$ws.Range("F21").Value = 1.002870546799215
$ws.Range("G21").Value = 0.0007711437311507585
$ws.Range("H21").Value = 0.03398707486947439
$ws.Range("I21").Value = 0.02118666147670556
$ws.Range("P21").Value = 0.691096490122753
$ws.Range("Q21").Value = 2.862923816794051

$ws.Range("B22").Value = 4.551661220263384
$ws.Range("C22").Value = 1.213745153735829
$ws.Range("D22").Value = 0.05114237283794409
$ws.Range("E22").Value = 2.398415709812426
$ws.Range("F22").Value = 1.070752524721385
$ws.Range("G22").Value = 0.000768271936691179
$ws.Range("H22").Value = 0.03907437156944749
$ws.Range("I22").Value = 0.02456729336757402
$ws.Range("P22").Value = 0.6855652115969377
$ws.Range("Q22").Value = 3.045853247365073

$ws.Range("B23").Value = 4.388043389984887
$ws.Range("C23").Value = 1.170230267484214
$ws.Range("D23").Value = 0.04958402825251795
$ws.Range("E23").Value = 2.309816508383619
$ws.Range("F23").Value = 1.03681460987049
$ws.Range("G23").Value = 0.0007697827481974752
$ws.Range("H23").Value = 0.03633694744138571
$ws.Range("I23").Value = 0.02250664103468569
$ws.Range("P23").Value = 0.6867202211229895
$ws.Range("Q23").Value = 2.95620379021554

$ws.Range("B24").Value = 3.759689213248521
$ws.Range("C24").Value = 1.012181686639536
$ws.Range("D24").Value = 0.04394261903344443
$ws.Range("E24").Value = 1.976461631646146
$ws.Range("F24").Value = 0.9067763043870798
$ws.Range("G24").Value = 0.0007756338485891714
$ws.Range("H24").Value = 0.02669420700939273
$ws.Range("I24").Value = 0.01574906392351672
$ws.Range("P24").Value = 0.6955573263841259
$ws.Range("Q24").Value = 2.610298883614547

$ws.Range("B25").Value = 3.082796153344248
$ws.Range("C25").Value = 0.8434117984034515
$ws.Range("D25").Value = 0.03782965901152835
$ws.Range("E25").Value = 1.620805740526038
$ws.Range("F25").Value = 0.7709755821200588
$ws.Range("G25").Value = 0.000782176244044853
$ws.Range("H25").Value = 0.01766107065060085
$ws.Range("I25").Value = 0.009773803154761396
$ws.Range("P25").Value = 0.7082226143536374
$ws.Range("Q25").Value = 2.250509085723536
